$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update F/G values for rows 56, 57 (Bolo Funcional items)
$ws.Range("F56").Value = 228
$ws.Range("F57").Value = 950

# Update F/G values for rows 68-72 (swap / adjust pairs)
$ws.Range("F68").Value = 210

$ws.Range("F69").Value = 146
$ws.Range("G69").Value = 100

$ws.Range("F70").Value = 120
$ws.Range("G70").Value = 160

$ws.Range("F71").Value = 118
$ws.Range("G71").Value = 160

$ws.Range("F72").Value = 118
$ws.Range("G72").Value = 160

# Update the view - scroll/selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("F58").Select()

$wb.Windows.Item(1).WindowState = -4143
